$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rows = @(
    @{ Row=37; Calidad="Especial"; Volumen=100; Min=6000; Max=7000; Prom=6500; PrecioKg=2167 },
    @{ Row=38; Calidad="Primera";  Volumen=120; Min=5000; Max=6000; Prom=5500; PrecioKg=1833 },
    @{ Row=39; Calidad="Segunda";  Volumen=140; Min=4000; Max=5000; Prom=4500; PrecioKg=1500 }
)

foreach ($r in $rows) {
    $row = $r.Row

    $ws.Cells.Item($row, 1).Value = 1
    $ws.Cells.Item($row, 2).Value = "Agrícola del Norte S.A. de Arica"
    $ws.Cells.Item($row, 3).Value = "Arica y Parinacota"

    $ws.Cells.Item($row, 4).Value = 44832
    $ws.Cells.Item($row, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"

    $ws.Cells.Item($row, 5).Value = 15
    $ws.Cells.Item($row, 6).Value = "Fruta"
    $ws.Cells.Item($row, 7).Value = 100101
    $ws.Cells.Item($row, 8).Value = "Berries"
    $ws.Cells.Item($row, 9).Value = 100112025
    $ws.Cells.Item($row, 10).Value = "Frutilla"
    $ws.Cells.Item($row, 11).Value = "Sin especificar"
    $ws.Cells.Item($row, 12).Value = $r.Calidad
    $ws.Cells.Item($row, 13).Value = $r.Volumen
    $ws.Cells.Item($row, 14).Value = $r.Min
    $ws.Cells.Item($row, 15).Value = $r.Max
    $ws.Cells.Item($row, 16).Value = $r.Prom
    $ws.Cells.Item($row, 17).Value = "$/bandeja 3 kilos"
    $ws.Cells.Item($row, 18).Value = "Región de Arica y Parinacota"
    $ws.Cells.Item($row, 19).Value = $r.PrecioKg
    $ws.Cells.Item($row, 20).Value = 3
}
